# Refresh the cryptos price list (GitHub Actions scheduled update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Note: some "Price" values are written with a leading apostrophe so Excel
# keeps them as text (matching the sheet's existing convention) instead of
# auto-coercing a numeric-looking string like "587.17" into a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.840.39"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.145.00"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'587.17"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'146.48"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.135.58"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'36.95"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.664.21"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "63.671.35"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "3.142.88"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'7.08"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "'464.32"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'12.91"
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D25").Value = "'81.12"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'9.33"
$ws.Range("E28").Value = "  +7.68%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.68"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'7.12"
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").Value = "'26.93"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'0.109"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "'436.35"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").Value = "'8.88"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "2.917.90"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D45").Value = "'0.280"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").Value = "'37.32"
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("D48").Value = "'126.85"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "'24.07"
$ws.Range("E51").Value = "  -3.46%  "
